$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.203.58'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '1.826.24'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.37'
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6008'
$ws.Range("E6").Value = '  -4.18%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07056'
$ws.Range("E8").Value = '  -5.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2790'
$ws.Range("E9").Value = '  -3.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.45'
$ws.Range("E10").Value = '  -5.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07636'
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = '1.825.94'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.781'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.000009923'
$ws.Range("E14").Value = '  -3.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6259'
$ws.Range("E15").Value = '  -7.52%  '
$ws.Range("D16").Value = '2.071.03'
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.99'
$ws.Range("E17").Value = '  -3.45%  '
$ws.Range("D18").Value = '29.185.57'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.827'
$ws.Range("E19").Value = '  -6.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '225.88'
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E22").Value = '  -5.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.994'
$ws.Range("E23").Value = '  -4.54%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.39'
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.004'
$ws.Range("E26").Value = '  -5.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1295'
$ws.Range("E27").Value = '  -4.17%  '
$ws.Range("E28").Value = '  -4.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.478'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06156'
$ws.Range("E30").Value = '  -15.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.440'
$ws.Range("E31").Value = '  -2.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.820'
$ws.Range("E32").Value = '  -5.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.793'
$ws.Range("E33").Value = '  -6.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.743'
$ws.Range("E35").Value = '  -4.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6397'
$ws.Range("E36").Value = '  -8.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.535'
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("D38").Value = '1.217.18'
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.734'
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("E40").Value = '  -5.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.529'
$ws.Range("E41").Value = '  -6.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9017'
$ws.Range("E42").Value = '  -4.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '1.984.37'
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.43'
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.57'
$ws.Range("E46").Value = '  -4.63%  '
$ws.Range("E47").Value = '  -1.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.487'
$ws.Range("E48").Value = '  -4.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.578'
$ws.Range("E49").Value = '  -8.55%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05509'
$ws.Range("E51").Value = '  -2.60%  '
